$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.345.99'
$ws.Range('E2').Value = '  -7.72%  '
$ws.Range('D3').Value = '2.889.33'
$ws.Range('E3').Value = '  -10.51%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '477.24'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -11.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '126.82'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.82%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '2.885.55'
$ws.Range('E8').Value = '  -10.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.405'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -11.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.68'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -12.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0976'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -14.94%  '
$ws.Range('E12').Value = '  -15.17%  '
$ws.Range('E13').Value = '  -3.70%  '
$ws.Range('D14').Value = '3.386.43'
$ws.Range('E14').Value = '  -10.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.86'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -12.01%  '
$ws.Range('D16').Value = '54.295.54'
$ws.Range('E16').Value = '  -7.89%  '
$ws.Range('D17').Value = '2.893.78'
$ws.Range('E17').Value = '  -10.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000136'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -14.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.24'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -11.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.66'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -12.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.13'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -13.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '310.96'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -14.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.450'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -13.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '59.81'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -15.28%  '
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.154'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -9.83%  '
$ws.Range('D29').Value = '0.0₃0826'
$ws.Range('E29').Value = '  -14.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.29'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -11.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.15'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.24'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -12.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.20'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -12.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.63'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -15.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.28'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -13.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '140.04'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -13.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.49'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -14.62%  '
$ws.Range('E38').Value = '  -15.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.03'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -12.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0622'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -12.15%  '
$ws.Range('D41').Value = '2.919.34'
$ws.Range('E41').Value = '  -10.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '35.46'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -13.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.965'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -12.55%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.602'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -15.94%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.44'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -14.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.33'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -11.95%  '
$ws.Range('D48').Value = '2.065.58'
$ws.Range('E48').Value = '  -10.23%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.36'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -15.14%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.07'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -13.47%  '
$ws.Range('E51').Value = '  -11.84%  '
